$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.79%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'49.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.72%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.321"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.33%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08067"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.52%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.607"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.38%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.353"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'28.76%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.643"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.77%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1280"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.91%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1972"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'6.04%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09627"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.72%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04719"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'13.53%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.27%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001327"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.63%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04209"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.52%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005858"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.97%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'2.469"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'5.94%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3508"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'4.36%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.167"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.76%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1381"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.23%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3092"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.71%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001295"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.46%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004294"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.52%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001351"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.72%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003537"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02732"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'8.77%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05958"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'12.01%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01080"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'93.08%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008041"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.41%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'7.31%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007557"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.07%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007898"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'5.49%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3220"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.91%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006981"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'4.39%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.59%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05541"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'27.41%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.37%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.59%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.59%"
$ws.Range("E51").Style = "Normal"
